$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1): descriptive Spanish headers -> short field names ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case the Spanish connector words ("de"/"del"/"la"/"las"/"el"/"los"/"y") ---
# --- in state/municipality name cells, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga" ---
$ws.Range("B8").Value = 'Pabellón De Arteaga'
$ws.Range("B9").Value = 'Rincón De Romos'
$ws.Range("B10").Value = 'San Francisco De Los Romo'
$ws.Range("B11").Value = 'San José De Gracia'
$ws.Range("B35").Value = 'Amatenango De La Frontera'
$ws.Range("B36").Value = 'Amatenango Del Valle'
$ws.Range("B39").Value = 'Bejucal De Ocampo'
$ws.Range("B46").Value = 'Chiapa De Corzo'
$ws.Range("B50").Value = 'Comitán De Domínguez'
$ws.Range("B68").Value = 'Mazapa De Madero'
$ws.Range("B70").Value = 'Montecristo De Guerrero'
$ws.Range("B73").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B82").Value = 'San Cristóbal De Las Casas'
$ws.Range("B118").Value = 'Hidalgo Del Parral'
$ws.Range("B128").Value = 'San Francisco Del Oro'
$ws.Range("B148").Value = 'San Juan De Sabinas'
$ws.Range("A162").Value = 'Ciudad De México'
$ws.Range("B166").Value = 'Cuajimalpa De Morelos'
$ws.Range("B180").Value = 'Coneto De Comonfort'
$ws.Range("B191").Value = 'Nombre De Dios'
$ws.Range("B195").Value = 'Pánuco De Coronado'
$ws.Range("B201").Value = 'San Juan Del Río'
$ws.Range("A210").Value = 'Estado De México'
$ws.Range("B210").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B213").Value = 'Almoloya De Alquisiras'
$ws.Range("B214").Value = 'Almoloya De Juárez'
$ws.Range("B215").Value = 'Almoloya Del Río'
$ws.Range("B220").Value = 'Atizapán De Zaragoza'
$ws.Range("B227").Value = 'Chapa De Mota'
$ws.Range("B230").Value = 'Coacalco De Berriozábal'
$ws.Range("B237").Value = 'Ecatepec De Morelos'
$ws.Range("B242").Value = 'Ixtapan De La Sal'
$ws.Range("B243").Value = 'Ixtapan Del Oro'
$ws.Range("B255").Value = 'Naucalpan De Juárez'
$ws.Range("B263").Value = 'San Antonio La Isla'
$ws.Range("B264").Value = 'San Felipe Del Progreso'
$ws.Range("B265").Value = 'San Martín De Las Pirámides'
$ws.Range("B267").Value = 'San Simón De Guerrero'
$ws.Range("B277").Value = 'Tenango Del Aire'
$ws.Range("B278").Value = 'Tenango Del Valle'
$ws.Range("B286").Value = 'Tlalnepantla De Baz'
$ws.Range("B292").Value = 'Valle De Bravo'
$ws.Range("B293").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B294").Value = 'Villa De Allende'
$ws.Range("B295").Value = 'Villa Del Carbón'
$ws.Range("B307").Value = 'San Miguel De Allende'
$ws.Range("B308").Value = 'Apaseo El Alto'
$ws.Range("B309").Value = 'Apaseo El Grande'
$ws.Range("B316").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B320").Value = 'Jaral Del Progreso'
$ws.Range("B328").Value = 'Purísima Del Rincón'
$ws.Range("B332").Value = 'San Diego De La Unión'
$ws.Range("B334").Value = 'San Francisco Del Rincón'
$ws.Range("B336").Value = 'San Luis De La Paz'
$ws.Range("B337").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B339").Value = 'Silao De La Victoria'
$ws.Range("B344").Value = 'Valle De Santiago'
$ws.Range("B350").Value = 'Acapulco De Juárez'
$ws.Range("B352").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B353").Value = 'Alcozauca De Guerrero'
$ws.Range("B356").Value = 'Atenango Del Río'
$ws.Range("B357").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B359").Value = 'Atoyac De Álvarez'
$ws.Range("B360").Value = 'Ayutla De Los Libres'
$ws.Range("B363").Value = 'Buenavista De Cuéllar'
$ws.Range("B364").Value = 'Chilapa De Álvarez'
$ws.Range("B365").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B366").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B371").Value = 'Coyuca De Benítez'
$ws.Range("B372").Value = 'Coyuca De Catalán'
$ws.Range("B375").Value = 'Cuetzala Del Progreso'
$ws.Range("B376").Value = 'Cutzamala De Pinzón'
$ws.Range("B382").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B383").Value = 'Iguala De La Independencia'
$ws.Range("B385").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B388").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B403").Value = 'Taxco De Alarcón'
$ws.Range("B405").Value = 'Técpan De Galeana'
$ws.Range("B407").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B409").Value = 'Tixtla De Guerrero'
$ws.Range("B412").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B413").Value = 'Tlapa De Comonfort'
$ws.Range("B423").Value = 'Agua Blanca De Iturbide'
$ws.Range("B429").Value = 'Atotonilco El Grande'
$ws.Range("B433").Value = 'Cuautepec De Hinojosa'
$ws.Range("B439").Value = 'Huejutla De Reyes'
$ws.Range("B442").Value = 'Jacala De Ledezma'
$ws.Range("B447").Value = 'Mixquiahuala De Juárez'
$ws.Range("B448").Value = 'Molango De Escamilla'
$ws.Range("B450").Value = 'Omitlán De Juárez'
$ws.Range("B451").Value = 'Pachuca De Soto'
$ws.Range("B458").Value = 'Santiago De Anaya'
$ws.Range("B462").Value = 'Tenango De Doria'
$ws.Range("B464").Value = 'Tepehuacán De Guerrero'
$ws.Range("B465").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B466").Value = 'Tezontepec De Aldama'
$ws.Range("B471").Value = 'Tula De Allende'
$ws.Range("B472").Value = 'Tulancingo De Bravo'
$ws.Range("B473").Value = 'Villa De Tezontepec'
$ws.Range("B475").Value = 'Zacualtipán De Ángeles'
$ws.Range("B476").Value = 'Zapotlán De Juárez'
$ws.Range("B481").Value = 'Acatlán De Juárez'
$ws.Range("B482").Value = 'Ahualulco De Mercado'
$ws.Range("B485").Value = 'Atotonilco El Alto'
$ws.Range("B487").Value = 'Autlán De Navarro'
$ws.Range("B496").Value = 'Concepción De Buenos Aires'
$ws.Range("B503").Value = 'Encarnación De Díaz'
$ws.Range("B508").Value = 'Huejuquilla El Alto'
$ws.Range("B509").Value = 'Ixtlahuacán Del Río'
$ws.Range("B513").Value = 'Jilotlán De Los Dolores'
$ws.Range("B518").Value = 'La Manzanilla De La Paz'
$ws.Range("B519").Value = 'Lagos De Moreno'
$ws.Range("B525").Value = 'Ojuelos De Jalisco'
$ws.Range("B531").Value = 'San Juan De Los Lagos'
$ws.Range("B534").Value = 'San Miguel El Alto'
$ws.Range("B535").Value = 'San Sebastián Del Oeste'
$ws.Range("B536").Value = 'Santa María De Los Ángeles'
$ws.Range("B539").Value = 'Talpa De Allende'
$ws.Range("B540").Value = 'Tamazula De Gordiano'
$ws.Range("B545").Value = 'Teocuitatlán De Corona'
$ws.Range("B546").Value = 'Tepatitlán De Morelos'
$ws.Range("B547").Value = 'Tizapán El Alto'
$ws.Range("B548").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B555").Value = 'Unión De San Antonio'
$ws.Range("B556").Value = 'Unión De Tula'
$ws.Range("B557").Value = 'Valle De Juárez'
$ws.Range("B562").Value = 'Yahualica De González Gallo'
$ws.Range("B563").Value = 'Zacoalco De Torres'
$ws.Range("B566").Value = 'Zapotlán Del Rey'
$ws.Range("B567").Value = 'Zapotlán El Grande'
$ws.Range("B589").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B591").Value = 'Cojumatlán De Régules'
$ws.Range("B654").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B677").Value = 'Coatlán Del Río'
$ws.Range("B689").Value = 'Puente De Ixtla'
$ws.Range("B694").Value = 'Tetela Del Volcán'
$ws.Range("B695").Value = 'Tlaltizapán De Zapata'
$ws.Range("B701").Value = 'Zacualpan De Amilpas'
$ws.Range("B704").Value = 'Amatlán De Cañas'
$ws.Range("B708").Value = 'Ixtlán Del Río'
$ws.Range("B713").Value = 'Santa María Del Oro'
$ws.Range("B730").Value = 'Mier Y Noriega'
$ws.Range("B734").Value = 'San Nicolás De Los Garza'
$ws.Range("B739").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B743").Value = 'Ayoquezco De Aldama'
$ws.Range("B745").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B747").Value = 'Constancia Del Rosario'
$ws.Range("B749").Value = 'Cuilápam De Guerrero'
$ws.Range("B750").Value = 'El Barrio De La Soledad'
$ws.Range("B751").Value = 'Fresnillo De Trujano'
$ws.Range("B753").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B754").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B755").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B756").Value = 'Huautla De Jiménez'
$ws.Range("B758").Value = 'Ixtlán De Juárez'
$ws.Range("B759").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B765").Value = 'Mariscala De Juárez'
$ws.Range("B766").Value = 'Mártires De Tacubaya'
$ws.Range("B768").Value = 'Mazatlán Villa De Flores'
$ws.Range("B770").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B772").Value = 'Oaxaca De Juárez'
$ws.Range("B773").Value = 'Ocotlán De Morelos'
$ws.Range("B774").Value = 'Pinotepa De Don Luis'
$ws.Range("B775").Value = 'Putla Villa De Guerrero'
$ws.Range("B776").Value = 'Rojas De Cuauhtémoc'
$ws.Range("B787").Value = 'San Antonino El Alto'
$ws.Range("B798").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B817").Value = 'San Juan Del Estado'
$ws.Range("B818").Value = 'San Juan Del Río'
$ws.Range("B848").Value = 'San Miguel Del Puerto'
$ws.Range("B855").Value = 'San Pablo Villa De Mitla'
$ws.Range("B858").Value = 'San Pedro El Alto'
$ws.Range("B870").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B883").Value = 'Santa Cruz De Bravo'
$ws.Range("B888").Value = 'Santa Inés Del Monte'
$ws.Range("B897").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B932").Value = 'Santo Domingo De Morelos'
$ws.Range("B942").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B943").Value = 'Tataltepec De Valdés'
$ws.Range("B944").Value = 'Teotitlán De Flores Magón'
$ws.Range("B945").Value = 'Teotitlán Del Valle'
$ws.Range("B946").Value = 'Tepelmeme Villa De Morelos'
$ws.Range("B948").Value = 'Tlacolula De Matamoros'
$ws.Range("B952").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B953").Value = 'Villa De Etla'
$ws.Range("B954").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B955").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B956").Value = 'Villa De Zaachila'
$ws.Range("B959").Value = 'Villa Sola De Vega'
$ws.Range("B960").Value = 'Villa Talea De Castro'
$ws.Range("B963").Value = 'Zimatlán De Álvarez'
$ws.Range("B980").Value = 'Ayotoxco De Guerrero'
$ws.Range("B983").Value = 'Chalchicomula De Sesma'
$ws.Range("B990").Value = 'Chila De La Sal'
$ws.Range("B1002").Value = 'Cuayuca De Andrade'
$ws.Range("B1003").Value = 'Cuetzalan Del Progreso'
$ws.Range("B1017").Value = 'Huehuetlán El Chico'
$ws.Range("B1018").Value = 'Huehuetlán El Grande'
$ws.Range("B1021").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B1024").Value = 'Izúcar De Matamoros'
$ws.Range("B1031").Value = 'Los Reyes De Juárez'
$ws.Range("B1040").Value = 'Palmar De Bravo'
$ws.Range("B1048").Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range("B1060").Value = 'San Salvador El Seco'
$ws.Range("B1061").Value = 'San Salvador El Verde'
$ws.Range("B1067").Value = 'Tecali De Herrera'
$ws.Range("B1075").Value = 'Tepanco De López'
$ws.Range("B1076").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1081").Value = 'Tepexi De Rodríguez'
$ws.Range("B1082").Value = 'Tetela De Ocampo'
$ws.Range("B1087").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1114").Value = 'Amealco De Bonfil'
$ws.Range("B1116").Value = 'Cadereyta De Montes'
$ws.Range("B1120").Value = 'Jalpan De Serra'
$ws.Range("B1121").Value = 'Landa De Matamoros'
$ws.Range("B1123").Value = 'Pinal De Amoles'
$ws.Range("B1126").Value = 'San Juan Del Río'
$ws.Range("B1135").Value = 'Axtla De Terrazas'
$ws.Range("B1141").Value = 'Ciudad Del Maíz'
$ws.Range("B1147").Value = 'Mexquitic De Carmona'
$ws.Range("B1152").Value = 'San Ciro De Acosta'
$ws.Range("B1156").Value = 'Santa María Del Río'
$ws.Range("B1158").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1167").Value = 'Villa De Arista'
$ws.Range("B1168").Value = 'Villa De Arriaga'
$ws.Range("B1169").Value = 'Villa De Ramos'
$ws.Range("B1170").Value = 'Villa De Reyes'
$ws.Range("B1201").Value = 'Nacozari De García'
$ws.Range("B1217").Value = 'Jalpa De Méndez'
$ws.Range("B1245").Value = 'Soto La Marina'
$ws.Range("B1258").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1264").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1266").Value = 'Muñoz De Domingo Arenas'
$ws.Range("B1269").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1270").Value = 'San Pablo Del Monte'
$ws.Range("B1273").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1276").Value = 'Tetla De La Solidaridad'
$ws.Range("B1296").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1300").Value = 'Amatlán De Los Reyes'
$ws.Range("B1307").Value = 'Boca Del Río'
$ws.Range("B1311").Value = 'Cazones De Herrera'
$ws.Range("B1325").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1326").Value = 'Cosautlán De Carvajal'
$ws.Range("B1340").Value = 'Hueyapan De Ocampo'
$ws.Range("B1341").Value = 'Ignacio De La Llave'
$ws.Range("B1344").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1345").Value = 'Ixhuatlán De Madero'
$ws.Range("B1346").Value = 'Ixhuatlán Del Café'
$ws.Range("B1355").Value = 'Juchique De Ferrer'
$ws.Range("B1358").Value = 'Landero Y Coss'
$ws.Range("B1361").Value = 'Las Vigas De Ramírez'
$ws.Range("B1362").Value = 'Lerdo De Tejada'
$ws.Range("B1366").Value = 'Martínez De La Torre'
$ws.Range("B1368").Value = 'Medellín De Bravo'
$ws.Range("B1382").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1386").Value = 'Paso De Ovejas'
$ws.Range("B1387").Value = 'Paso Del Macho'
$ws.Range("B1390").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1396").Value = 'Sayula De Alemán'
$ws.Range("B1397").Value = 'Soledad De Doblado'
$ws.Range("B1418").Value = 'Tlacotepec De Mejía'
$ws.Range("B1428").Value = 'Vega De Alatorre'
$ws.Range("B1436").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1451").Value = 'Concepción Del Oro'
$ws.Range("B1470").Value = 'Moyahua De Estrada'
$ws.Range("B1471").Value = 'Nochistlán De Mejía'
$ws.Range("B1472").Value = 'Noria De Ángeles'
$ws.Range("B1483").Value = 'Teúl De González Ortega'
$ws.Range("B1484").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1488").Value = 'Villa De Cos'

# --- Tiny floating point recompute delta on D350 (Guerrero / Acapulco de Juarez pct) ---
$ws.Range("D350").Value = 0.009100733746658323

# --- Drop the trailing footnote/metadata rows (1497-1501); row 1496 already blank ---
$ws.Rows("1497:1501").Delete()
